# Item.xlsx — "add Item Config And So on"
# Adds a new "Icon" row to the Property sheet's item-field table, plus the
# (hidden) LOCAL_MYSQL_DATE_FORMAT defined name that shows up alongside it,
# and leaves the active selection where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Workbook-level defined name (hidden) -------------------------------
$n = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", "=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&"" ""&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)")
$n.Visible = $false

# --- New row 19: Icon field definition ----------------------------------
$ws.Range("A19").Value = "Icon"
$ws.Range("B19").Value = "string"
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = $false
$ws.Range("E19").Value = $false
$ws.Range("F19").Value = $true
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = "Friend"
$ws.Range("J19").Value = "物品显示Icon"

# match the text-formatted style used by the other rows' Id/Type/RelationValue/Desc cells
$ws.Range("A19").NumberFormat = "@"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("I19").NumberFormat = "@"
$ws.Range("J19").NumberFormat = "@"

# --- Leave the selection where the new row was entered ------------------
[void]$ws.Range("J17").Select()
